# Update the 'scraped_at' timestamps (column K) on the 'snapshot' sheet
# to reflect the newest scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$updates = @(
    @{ Row = 2; Value = "2025-11-23T11:28:38.156335+00:00" },
    @{ Row = 3; Value = "2025-11-23T11:28:38.156368+00:00" },
    @{ Row = 4; Value = "2025-11-23T11:28:40.211442+00:00" },
    @{ Row = 5; Value = "2025-11-23T11:28:40.211472+00:00" },
    @{ Row = 6; Value = "2025-11-23T11:28:42.262328+00:00" },
    @{ Row = 7; Value = "2025-11-23T11:28:44.274086+00:00" },
    @{ Row = 8; Value = "2025-11-23T11:28:46.365788+00:00" },
    @{ Row = 9; Value = "2025-11-23T11:28:46.365805+00:00" },
    @{ Row = 10; Value = "2025-11-23T11:28:46.365813+00:00" },
    @{ Row = 11; Value = "2025-11-23T11:28:48.410769+00:00" },
    @{ Row = 12; Value = "2025-11-23T11:28:50.511452+00:00" },
    @{ Row = 13; Value = "2025-11-23T11:28:53.030961+00:00" },
    @{ Row = 14; Value = "2025-11-23T11:28:55.106932+00:00" },
    @{ Row = 15; Value = "2025-11-23T11:28:57.241116+00:00" },
    @{ Row = 16; Value = "2025-11-23T11:29:01.771944+00:00" },
    @{ Row = 17; Value = "2025-11-23T11:29:01.771975+00:00" },
    @{ Row = 18; Value = "2025-11-23T11:29:04.252782+00:00" },
    @{ Row = 19; Value = "2025-11-23T11:29:04.252825+00:00" },
    @{ Row = 20; Value = "2025-11-23T11:29:04.252836+00:00" },
    @{ Row = 21; Value = "2025-11-23T11:29:06.388260+00:00" },
    @{ Row = 22; Value = "2025-11-23T11:29:06.388277+00:00" },
    @{ Row = 23; Value = "2025-11-23T11:29:08.458157+00:00" },
    @{ Row = 24; Value = "2025-11-23T11:29:08.458185+00:00" },
    @{ Row = 25; Value = "2025-11-23T11:29:08.458194+00:00" },
    @{ Row = 26; Value = "2025-11-23T11:29:08.458202+00:00" },
    @{ Row = 27; Value = "2025-11-23T11:29:10.581664+00:00" },
    @{ Row = 28; Value = "2025-11-23T11:29:10.581708+00:00" },
    @{ Row = 29; Value = "2025-11-23T11:29:13.075546+00:00" },
    @{ Row = 30; Value = "2025-11-23T11:29:13.075573+00:00" },
    @{ Row = 31; Value = "2025-11-23T11:29:13.075589+00:00" },
    @{ Row = 32; Value = "2025-11-23T11:29:13.075603+00:00" },
    @{ Row = 33; Value = "2025-11-23T11:29:15.622114+00:00" },
    @{ Row = 34; Value = "2025-11-23T11:29:15.622141+00:00" },
    @{ Row = 35; Value = "2025-11-23T11:29:19.887165+00:00" },
    @{ Row = 36; Value = "2025-11-23T11:29:19.887180+00:00" },
    @{ Row = 37; Value = "2025-11-23T11:29:22.706590+00:00" },
    @{ Row = 38; Value = "2025-11-23T11:29:22.706617+00:00" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 11).Value = $u.Value
}
